$d = $word.ActiveDocument

$pairs = @(
    @("2024-10-28 Monday", "2024-10-29 Tuesday"),
    @("828÷3=276, 0", "615÷6=102, 3"),
    @("401÷8=50, 1", "226÷2=113, 0"),
    @("191÷7=27, 2", "976÷3=325, 1"),
    @("296÷8=37, 0", "720÷7=102, 6"),
    @("277÷4=69, 1", "657÷4=164, 1"),
    @("223÷5=44, 3", "432÷5=86, 2"),
    @("123÷3=41, 0", "347÷7=49, 4"),
    @("564÷9=62, 6", "762÷3=254, 0"),
    @("320÷8=40, 0", "219÷6=36, 3"),
    @("662÷9=73, 5", "531÷7=75, 6"),
    @("586÷4=146, 2", "841÷8=105, 1"),
    @("506÷3=168, 2", "612÷3=204, 0"),
    @("833÷6=138, 5", "876÷4=219, 0"),
    @("887÷4=221, 3", "835÷9=92, 7"),
    @("740÷6=123, 2", "986÷4=246, 2"),
    @("822÷4=205, 2", "951÷2=475, 1"),
    @("917÷4=229, 1", "626÷4=156, 2"),
    @("324÷2=162, 0", "364÷3=121, 1"),
    @("231÷4=57, 3", "319÷9=35, 4"),
    @("995÷5=199, 0", "170÷4=42, 2"),
    @("415÷6=69, 1", "119÷7=17, 0"),
    @("515÷3=171, 2", "845÷9=93, 8"),
    @("470÷9=52, 2", "220÷7=31, 3"),
    @("844÷6=140, 4", "488÷5=97, 3"),
    @("354÷6=59, 0", "760÷6=126, 4")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "done"
